# The underlying edit swaps the contents of ppt/theme/theme1.xml (the
# slide master's theme, originally the "Integral" theme) and
# ppt/theme/theme2.xml (the notes master's theme, originally the
# default "Office Theme"). In practice the fontScheme/fmtScheme of the
# two theme parts are already byte-identical, so the only observable
# difference between "Integral" and "Office Theme" is the 12-colour
# clrScheme. Re-apply the "Office Theme" colour scheme to the
# presentation's (slide-master) theme via the ThemeColorScheme object,
# which is the supported automation surface for editing a theme's
# clrScheme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
